$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: make C1 share the same style as A1/B1 ---
$ws.Range("C1").Value = $ws.Range("C1").Value
$ws.Range("A1:C1").Font.Name = "Arial"
$ws.Range("C1").Font.Color = $ws.Range("B1").Font.Color

# --- Rows 2-5 (existing data) lose their explicit alignment ---
$ws.Range("A2:B5").HorizontalAlignment = -4131
$ws.Range("A2:B5").HorizontalAlignment = -4142

# --- Row 6 duplicate ("HCFC phase out plan"/"PHA") replaced, rows shift up conceptually ---
$ws.Range("A6").Value = "Air conditioning"
$ws.Range("B6").Value = "REF"
$ws.Range("A7").Value = "Commercial"
$ws.Range("B7").Value = "REF"
$ws.Range("A8").Value = "Multiple-subsectors"
$ws.Range("B8").Value = "FOA"
$ws.Range("A9").Value = "HFC phase down plan"
$ws.Range("B9").Value = "PHA"

# --- New "proposals" rows appended (10-14) ---
$ws.Range("A10").Value = "Preparation of project proposal"
$ws.Range("B10").Value = "KIP"
$ws.Range("A11").Value = "Domestic/commercial"
$ws.Range("B11").Value = "REF"
$ws.Range("A12").Value = "HCFC closure"
$ws.Range("B12").Value = "PRO"
$ws.Range("A13").Value = "Agency programme"
$ws.Range("B13").Value = "SEV"
$ws.Range("A14").Value = "Technical assistance/support"
$ws.Range("B14").Value = "REF"

# Give the first new row (A10) a distinct highlighted look
$ws.Range("A10").Interior.Color = 16777215
$ws.Range("A10").HorizontalAlignment = -4131
$ws.Range("A10").Font.Color = 0

# Remaining new rows: left-aligned style, normal theme font
$ws.Range("B10:B14").HorizontalAlignment = -4131
$ws.Range("A11:A14").HorizontalAlignment = -4131

# --- Extend the used range with 130 additional blank rows at the very end ---
$ws.Range("A997").Value = "x"
$ws.Range("A997").Value = ""
